$d = $word.ActiveDocument

function Rename-LogoShape($shape) {
    $alt = $shape.AlternativeText
    if ($alt -eq "BTec_Logo-Orange") {
        $shape.Name = "image2.jpg"
    } elseif ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shape.Name = "image1.png"
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($hi = 1; $hi -le 3; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                Rename-LogoShape $shapes.Item($i)
            }
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                Rename-LogoShape $shapes.Item($i)
            }
        }
    }
}
